$wb = $excel.ActiveWorkbook

# --- Sheet 1: "RO Non-Availability Dates" ---
# Insert a new column before F for the new Regional Office "RO05" /
# "White River Junction, VT" and give it the same formatting/width as its
# neighboring column.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F:F").Insert()
$ws1.Columns.Item(6).ColumnWidth = $ws1.Columns.Item(7).ColumnWidth

$ws1.Range("F2").Value = "RO05"
$ws1.Range("F3").Value = "White River Junction, VT"

# The non-availability dates for the new RO match the ones already present
# on the neighboring columns for the first few rows of data.
$ws1.Range("E4:E7").Copy()
$ws1.Range("F4:F7").PasteSpecial()

# --- Sheet 3: "RO Allocations" ---
# Insert a new row for RO05 / White River Junction, VT (keeping the table
# sorted by RO number, between RO04 and RO06).
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("8:8").Insert()
$ws3.Range("B8").Value = "White River Junction, VT"
$ws3.Range("C8").Value = "RO05"
$ws3.Range("D8").Value = 4
$ws3.Range("E8").Value = 0
$ws3.Range("F8").Value = 0
$ws3.Range("G8").Value = 0
$ws3.Range("H8").Value = 0

# --- Window / selection state ---
# The active sheet moves from "RO Allocations" to "RO Non-Availability
# Dates", with the newly added cell selected there, while "RO Allocations"
# keeps a plain single-cell selection.
$ws3.Range("C13").Select()
$ws1.Activate()
$ws1.Range("F3").Select()
